$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 4 through 8 (CBS News stays removed along with Ceska Televize,
# Ceske Radiokomunikace, Cetin, Channel 4 (UK))
$ws.Rows("4:8").Delete()

# Reorder remaining rows: A1=CBS (unchanged), A2=Cellnex, A3=CBS Studios International
$ws.Range("A2").Value = "Cellnex"
$ws.Range("A3").Value = "CBS Studios International"

# Update selection to match target state
$ws.Range("A2").Select()
